$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Narrow the one-person-template table: 9576 dxa (478.8pt) -> 8540 dxa
#    (427pt), and give the single row an explicit height of 4144 twips
#    (207.2pt), matching the author's commit (SDEIS display-width fix).
# ---------------------------------------------------------------------------
$t = $d.Tables.Item(1)

$newWidthPts = 8540 / 20      # 427
$newHeightPts = 4144 / 20     # 207.2

$t.PreferredWidth = $newWidthPts
$t.Columns.Item(1).Width = $newWidthPts
$t.Cell(1, 1).Width = $newWidthPts

# Row height only (no HeightRule touch) -> serializes as <w:trHeight w:val="4144"/>
$t.Rows.Item(1).Height = $newHeightPts

# ---------------------------------------------------------------------------
# 2. Wrap the "RoomNumber" run in spell-check proofing marks
#    (<w:proofErr w:type="spellStart"/> ... <w:proofErr w:type="spellEnd"/>).
#    Locate it with Find (Paragraphs.Item(n) indexing can go stale right
#    after the table resize above), then replace the whole paragraph
#    (text + its end-of-paragraph mark) in place via InsertXML so the
#    proofErr tags land as siblings of the run instead of splitting the
#    paragraph.
# ---------------------------------------------------------------------------
$findRange = $d.Content
$found = $findRange.Find.Execute("RoomNumber", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "RoomNumber paragraph not found"
}

$fullPara = $d.Range($findRange.Start, $findRange.End + 1)

$pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData>' +
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:body>' +
       '<w:p>' +
       '<w:proofErr w:type="spellStart"/>' +
       '<w:r><w:t>RoomNumber</w:t></w:r>' +
       '<w:proofErr w:type="spellEnd"/>' +
       '</w:p>' +
       '</w:body></w:document>' +
       '</pkg:xmlData></pkg:part></pkg:package>'

$fullPara.InsertXML($pkg)
